# Fix formatting on fastq purpose column: "fullRNASEQ" -> "fullRNASeq"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
